# Daily attendance processing - 2025-12-14 09:28:01
#
# Normalises the "Recorded By" column (G): for the specific recorded-by
# combinations produced by the automated session importer, the literal
# "System" marker is moved from the front of the comma-separated list to
# the back, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Only the exact known combinations are touched; anything already in the
# other order (or combinations such as "System, admin@admin.com") is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
